# Generate Report for Handoff
# Refresh the "Latest Handoff" timestamps for file
# 52abe100-60f7-4f8b-93a5-08fb336a0d6e (row 4 on every sheet) after a new
# handoff package was generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the
# 52abe100-... row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2017-02-21 02:25:55"

# zh-cn sheet: "Latest Handoff Datetime" (column H) for the same file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2017-02-21 02:25:40"

# de-de sheet: "Latest Handoff Datetime" (column H) for the same file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2017-02-21 02:25:55"
